# Insert a new data row at row 17 (pushing the existing rows 17-38 down to
# 18-39) and populate it with a new "Haba" price record for
# Terminal La Palmera de La Serena.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 17..38 down to 18..39, leaving row 17 free for the new record.
$ws.Rows.Item(17).Insert()

$ws.Range("A17").Value = 8
$ws.Range("B17").Value = "Terminal La Palmera de La Serena"
$ws.Range("C17").Value = "Coquimbo"
$ws.Range("D17").Value = 44874
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 100112026
$ws.Range("G17").Value = "Haba"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 7000
$ws.Range("M17").Value = 6500
$ws.Range("N17").Value = "`$/saco 25 kilos"
$ws.Range("O17").Value = "Provincia del Elquí"
$ws.Range("P17").Value = 260
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = "Hortaliza"
